# edit.ps1 -- apply "Updated notes on terms and concepts" changes
#
# Summary of the change (derived from the target XML diff):
#  1. "Policy Paradox" paragraph: "posited" -> "espoused", and the
#     hidden `_GoBack` bookmark (originally a collapsed bookmark right
#     after "Administrative Evil") is relocated to sit right after the
#     word "espoused" in this paragraph.
#  2. Several paragraphs that had been split into many small runs
#     (left over from tracked editing) get their runs consolidated back
#     together -- "The Benevolent Community", "Mob at the Gates", and
#     "Triumphant Individual". Run boundaries at `w:proofErr` markers,
#     and at boundaries between genuinely different edit sessions, are
#     preserved; only runs that belong to the same contiguous span of
#     text collapse into one run.
#
# The Word engine backing this script only restructures/merges runs
# when the replaced text content actually changes, so each merge below
# is done in two passes: first the target span is swapped for a unique
# placeholder token (forcing Word to rebuild that span as a single new
# run), then the placeholder is swapped back for the real final text.

$d = $word.ActiveDocument

function Merge-Span {
    param($OldText, $NewText, $TokenSuffix)

    $f = $d.Content.Find
    $f.ClearFormatting()
    $f.Text = $OldText
    $f.MatchWildcards = $false
    $found = $f.Execute()
    if (-not $found -or -not $f.Found) {
        throw "Merge-Span: could not find text: $OldText"
    }

    $token = "TOKEN_MERGE_" + $TokenSuffix
    $rng = $f.Parent
    $rng.Text = $token

    $f2 = $d.Content.Find
    $f2.ClearFormatting()
    $f2.Text = $token
    $found2 = $f2.Execute()
    if (-not $found2 -or -not $f2.Found) {
        throw "Merge-Span: could not find placeholder token for: $OldText"
    }
    $f2.Parent.Text = $NewText
}

$VT = [char]11     # manual line break (<w:br/>) as seen through Range.Text
$RQ = [char]8217    # right single quotation mark (’)

# ---------------------------------------------------------------------
# 1) "The Benevolent Community"
# ---------------------------------------------------------------------

$oldA1 = $VT + "One of four morality tales in Robert Reich" + $RQ + "s conceptualization of the policy process.  In this morality tale, members of the community take care of one another.  In the conservative interpretation, this generosity is provided only to the deserving.  In the liberal interpretation, this generosity is provided because everyone "
Merge-Span $oldA1 $oldA1 "A1"

$oldA2 = " simply by being members of the community."
Merge-Span $oldA2 $oldA2 "A2"

# ---------------------------------------------------------------------
# 2) "Mob at the Gates"
# ---------------------------------------------------------------------

$oldB1 = $VT + "One of four morality tales in Robert Reich" + $RQ + "s conceptualization of the policy process.  In this morality tale, hordes of people want to enter our society.  In the conservative interpretation, these hordes are a threat to our way of life that must be stopped.  In the liberal interpretation,"
Merge-Span $oldB1 $oldB1 "B1"

# ---------------------------------------------------------------------
# 3) "Triumphant Individual"
# ---------------------------------------------------------------------

$oldC1 = $VT + "One of four morality tales in Robert Reich" + $RQ + "s conceptualization of the policy process.  In this morality tale, anybody can achieve anything in our society.  In the conservative interpretation, those that accomplish great things do so through their own resolve and determination.  In the liberal interpretation, those that accomplish great things do so with significant help from community resources and support."
Merge-Span $oldC1 $oldC1 "C1"

# ---------------------------------------------------------------------
# 4) "Policy Paradox": posited -> espoused
# ---------------------------------------------------------------------

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "posited"
$find.Replacement.Text = "espoused"
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, "espoused", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Move the hidden `_GoBack` bookmark so it sits right after
#    "espoused" in the "Policy Paradox" paragraph (collapsed bookmark).
#    Re-adding a bookmark with an existing name relocates it.
# ---------------------------------------------------------------------

$fEsp = $d.Content.Find
$fEsp.ClearFormatting()
$fEsp.Text = "espoused"
$fEsp.Execute() | Out-Null
$afterEspoused = $fEsp.Parent.End
$bmRange = $d.Range($afterEspoused, $afterEspoused)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
